$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 433.33334
$ws.Range("I42").Value = 300
$ws.Range("J42").Value = 500
$ws.Range("K42").Value = 900
$ws.Range("L42").Value = 1500
$ws.Range("M42").Value = -670
$ws.Range("N42").Value = -1960

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1374.5671
$ws.Range("I132").Value = 797.06665
$ws.Range("K132").Value = 2391.19995
$ws.Range("M132").Value = 138.8000499999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 44090.7
$ws.Range("J133").Value = 44090.7
$ws.Range("L133").Value = 44090.7
$ws.Range("N133").Value = -54210.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2210.0688
$ws.Range("I137").Value = 1427
$ws.Range("J137").Value = 8996.666999999999
$ws.Range("K137").Value = 4281
$ws.Range("L137").Value = 26990.001
$ws.Range("M137").Value = -1731
$ws.Range("N137").Value = -32090.001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2289.3
$ws.Range("I138").Value = 1137.9286
$ws.Range("J138").Value = 3123.0518
$ws.Range("K138").Value = 3413.7858
$ws.Range("L138").Value = 9369.1554
$ws.Range("M138").Value = 1726.2142
$ws.Range("N138").Value = -19649.1554

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1939.0975
$ws.Range("I141").Value = 1254.4546
$ws.Range("J141").Value = 4763.25
$ws.Range("K141").Value = 3763.3638
$ws.Range("L141").Value = 14289.75
$ws.Range("M141").Value = 1416.6362
$ws.Range("N141").Value = -24649.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 329067.44
$ws.Range("I61").Value = 8768.764999999999
$ws.Range("J61").Value = 718001.5600000001
$ws.Range("K61").Value = 8768.764999999999
$ws.Range("L61").Value = 718001.5600000001
$ws.Range("M61").Value = -8556.764999999999
$ws.Range("N61").Value = -718425.5600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2034.1923
$ws.Range("I74").Value = 1773.3334
$ws.Range("J74").Value = 2389.9092
$ws.Range("K74").Value = 1773.3334
$ws.Range("L74").Value = 2389.9092
$ws.Range("M74").Value = -899.3334
$ws.Range("N74").Value = -4137.9092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2034.1923
$ws.Range("I77").Value = 1773.3334
$ws.Range("J77").Value = 2389.9092
$ws.Range("K77").Value = 8866.666999999999
$ws.Range("L77").Value = 11949.546
$ws.Range("M77").Value = -4498.666999999999
$ws.Range("N77").Value = -20685.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4346.8
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 4808.5
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 4808.5
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -5620.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 4346.8
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 4808.5
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 4808.5
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -7616.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3128657.8
$ws.Range("I132").Value = 2096.2778
$ws.Range("J132").Value = 7148522.5
$ws.Range("K132").Value = 6288.8334
$ws.Range("L132").Value = 21445567.5
$ws.Range("M132").Value = -3758.8334
$ws.Range("N132").Value = -21450627.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 329067.44
$ws.Range("I136").Value = 8768.764999999999
$ws.Range("J136").Value = 718001.5600000001
$ws.Range("K136").Value = 26306.295
$ws.Range("L136").Value = 2154004.68
$ws.Range("M136").Value = -23756.295
$ws.Range("N136").Value = -2159104.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1937.25
$ws.Range("I86").Value = 1928.2858
$ws.Range("K86").Value = 1928.2858
$ws.Range("M86").Value = -805.2858000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1937.25
$ws.Range("I89").Value = 1928.2858
$ws.Range("K89").Value = 9641.429
$ws.Range("M89").Value = -4025.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 678626.4
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 678626.4
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 678626.4
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -679216.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 678626.4
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 678626.4
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 678626.4
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -679030.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 239750.36
$ws.Range("I58").Value = 1421.9615
$ws.Range("K58").Value = 1421.9615
$ws.Range("M58").Value = -1218.9615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 374389.94
$ws.Range("I134").Value = 4420.8423
$ws.Range("J134").Value = 1253066.5
$ws.Range("K134").Value = 13262.5269
$ws.Range("L134").Value = 3759199.5
$ws.Range("M134").Value = -10727.5269
$ws.Range("N134").Value = -3764269.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 239750.36
$ws.Range("I136").Value = 1421.9615
$ws.Range("K136").Value = 4265.8845
$ws.Range("M136").Value = -1715.8845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1613971.8
$ws.Range("I131").Value = 4167099.5
$ws.Range("J131").Value = 1470.0264
$ws.Range("K131").Value = 12501298.5
$ws.Range("L131").Value = 4410.0792
$ws.Range("M131").Value = -12496258.5
$ws.Range("N131").Value = -14490.0792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8350.3125
$ws.Range("I80").Value = 9236.071
$ws.Range("J80").Value = 2150
$ws.Range("K80").Value = 9236.071
$ws.Range("L80").Value = 2150
$ws.Range("M80").Value = -8238.071
$ws.Range("N80").Value = -4146

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 8350.3125
$ws.Range("I83").Value = 9236.071
$ws.Range("J83").Value = 2150
$ws.Range("K83").Value = 46180.355
$ws.Range("L83").Value = 10750
$ws.Range("M83").Value = -41188.355
$ws.Range("N83").Value = -20734

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3955.85
$ws.Range("I132").Value = 3497.9
$ws.Range("J132").Value = 5329.7
$ws.Range("K132").Value = 10493.7
$ws.Range("L132").Value = 15989.1
$ws.Range("M132").Value = -7963.700000000001
$ws.Range("N132").Value = -21049.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 52634920
$ws.Range("I40").Value = 58826796
$ws.Range("K40").Value = 58826796
$ws.Range("M40").Value = -58826660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2719639
$ws.Range("I122").Value = 5109433.5
$ws.Range("J122").Value = 628568.75
$ws.Range("K122").Value = 15328300.5
$ws.Range("L122").Value = 1885706.25
$ws.Range("M122").Value = -15325850.5
$ws.Range("N122").Value = -1890606.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 70715
$ws.Range("J127").Value = 70715
$ws.Range("L127").Value = 70715
$ws.Range("N127").Value = -80635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9998.1875
$ws.Range("I136").Value = 6879.72
$ws.Range("J136").Value = 21135.572
$ws.Range("K136").Value = 20639.16
$ws.Range("L136").Value = 63406.716
$ws.Range("M136").Value = -18089.16
$ws.Range("N136").Value = -68506.716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2105.5715
$ws.Range("I122").Value = 2123.1667
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6369.500100000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3919.500100000001
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2579.5
$ws.Range("I136").Value = 2399.4
$ws.Range("J136").Value = 2812.963
$ws.Range("K136").Value = 7198.200000000001
$ws.Range("L136").Value = 8438.889000000001
$ws.Range("M136").Value = -4648.200000000001
$ws.Range("N136").Value = -13538.889
